$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark rows 6-12 (ids 10000004..10000010) as Available = TRUE
$ws.Range("B6:B12").Value = $true

# Update row 30 (id 10000028): rename Name, clear Description
$ws.Range("C30").Value = "Part 29"
$ws.Range("D30").Value = ""

# Update row 31 (id 10000029): rename Name, clear Description
$ws.Range("C31").Value = "Part 10xx"
$ws.Range("D31").Value = ""

# Update row 32 (id 10000030): rename Name
$ws.Range("C32").Value = "art 10xx"

# Delete the old row 33 (id 10000032, Name=200) entirely - it is removed from the sheet,
# shifting row 34 (id 10000031, Name="Teil 281") up into row 33.
$ws.Rows(33).Delete()

# Update the (now) row 33 (id 10000031): rename Name to "Oi"
$ws.Range("C33").Value = "Oi"
